$wb = $excel.ActiveWorkbook

$hdrSheet  = $wb.Worksheets.Item("Batch_Header")
$detSheet  = $wb.Worksheets.Item("Batch_Detail")
$miscSheet = $wb.Worksheets.Item("Batch_Miscellaneous")

$oldBatchId = 250080
$newBatchId = 227799

# --- Batch_Header (sheet1): BATCHID (col B) + HDRDEBIT/HDRCREDIT (cols BO/BP) ---
$hdrDebit  = @(769037.95, 788794.53, 601426.4299999999, 104316.03)
$hdrCredit = @(-553147.46, -2898.56, -275606.79, -65195.99)

for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 2
    $hdrSheet.Range("B$row").Value = $newBatchId
    $hdrSheet.Range("BO$row").Value = $hdrDebit[$i]
    $hdrSheet.Range("BP$row").Value = $hdrCredit[$i]
}

# --- Batch_Detail (sheet2): header label CO1 = VALUES, BATCHID (col C), CO column -> 0 ---
$detSheet.Range("CO1").Value = "VALUES"

$lastRow = 97
for ($row = 2; $row -le $lastRow; $row++) {
    $detSheet.Range("C$row").Value = $newBatchId
    $detSheet.Range("CO$row").Value = 0
}

# --- Batch_Miscellaneous (sheet3): BATCHID (col A) ---
for ($row = 2; $row -le 5; $row++) {
    $miscSheet.Range("A$row").Value = $newBatchId
}
